$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.913.38'
$ws.Range('E2').Value = '  -2.08%  '
$ws.Range('D3').Value = '1.898.38'
$ws.Range('E3').Value = '  -4.01%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '324.22'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.08%  '
$ws.Range('E6').Value = '  -0.24%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4584'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -1.72%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3811'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.69%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07709'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -3.19%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.9745'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.95%  '
$ws.Range('E11').Value = '  -3.94%  '
$ws.Range('D12').Value = '1.887.19'
$ws.Range('E12').Value = '  -4.88%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.919'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -3.97%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.629'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -3.83%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.07017'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.06%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.004'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.22%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '83.49'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -4.93%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000009454'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -4.99%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '16.59'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -4.22%  '
$ws.Range('E20').Value = '  -0.26%  '
$ws.Range('D21').Value = '28.889.17'
$ws.Range('E21').Value = '  -2.18%  '
$ws.Range('E22').Value = '  -4.95%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.83'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.28%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.094'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.68%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '157.79'
$ws.Range('D25').Style = 'Normal'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '18.99'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.83%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '5.604'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.35%  '
$ws.Range('E28').Value = '  -2.04%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.829'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -4.23%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.09228'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.14%  '
$ws.Range('E31').Value = '  -4.24%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.073'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.22%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.992'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -6.19%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.05662'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.77%  '
$ws.Range('E36').Value = '  -2.85%  '
$ws.Range('E37').Value = '  -0.17%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02031'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -3.55%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.5464'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -4.60%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '7.363'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -5.47%  '
$ws.Range('E41').Value = '  -3.10%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '9.240'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -4.40%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.756'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.21%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.5141'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -4.21%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '11.19'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -5.64%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.06808'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.77%  '
$ws.Range('B47').Value = 'PEPE'
$ws.Range('C47').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.000002612'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -16.67%  '
$ws.Range('B48').Value = 'RenderToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.061'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -5.55%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '109.83'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -3.93%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.766'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -3.64%  '
$ws.Range('E51').Value = '  -0.19%  '
